$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.677.64'
$ws.Range("E2").Value = '  +1.61%  '
$ws.Range("D3").Value = '1.603.01'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.37'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.516'
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.05'
$ws.Range("E8").Value = '  +6.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0603'
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '1.832.37'
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("D13").Value = '1.605.26'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.548'
$ws.Range("E14").Value = '  +4.34%  '
$ws.Range("D15").Value = '29.676.66'
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.94'
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.56'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.81'
$ws.Range("E19").Value = '  +4.30%  '
$ws.Range("D20").Value = '0.0₃0698'
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.44'
$ws.Range("E23").Value = '  +2.41%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.36'
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("E26").Value = '  +2.04%  '
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.44'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("E30").Value = '  +2.45%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("D34").Value = '1.427.25'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.93'
$ws.Range("E35").Value = '  +3.79%  '
$ws.Range("E36").Value = '  +3.52%  '
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("E39").Value = '  +2.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '58.43'
$ws.Range("E40").Value = '  +7.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.546'
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("E42").Value = '  +5.84%  '
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.818'
$ws.Range("E44").Value = '  +3.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  +2.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.971'
$ws.Range("E47").Value = '  +15.63%  '
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").Value = '1.742.46'
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.76'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0525'
$ws.Range("E51").Value = '  +1.61%  '
